# Apply cryptos-list price/volume refresh (GitHub Actions data update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Ref = "D2"; Val = "56.411.79"},
    @{Ref = "E2"; Val = "  -2.63%  "},
    @{Ref = "D3"; Val = "2.944.89"},
    @{Ref = "E3"; Val = "  -3.94%  "},
    @{Ref = "E4"; Val = "  +0.21%  "},
    @{Ref = "D5"; Val = "492.87"},
    @{Ref = "E5"; Val = "  -6.58%  "},
    @{Ref = "D6"; Val = "134.23"},
    @{Ref = "E6"; Val = "  -6.64%  "},
    @{Ref = "E7"; Val = "  +0.49%  "},
    @{Ref = "D8"; Val = "0.423"},
    @{Ref = "E8"; Val = "  -5.73%  "},
    @{Ref = "D9"; Val = "7.13"},
    @{Ref = "E9"; Val = "  -7.09%  "},
    @{Ref = "D10"; Val = "0.105"},
    @{Ref = "E10"; Val = "  -7.29%  "},
    @{Ref = "D11"; Val = "0.349"},
    @{Ref = "E11"; Val = "  -5.94%  "},
    @{Ref = "D12"; Val = "3.474.55"},
    @{Ref = "E12"; Val = "  -3.31%  "},
    @{Ref = "E13"; Val = "  -3.12%  "},
    @{Ref = "D14"; Val = "25.53"},
    @{Ref = "E14"; Val = "  -7.14%  "},
    @{Ref = "E15"; Val = "  -9.21%  "},
    @{Ref = "D16"; Val = "56.663.60"},
    @{Ref = "E16"; Val = "  -2.23%  "},
    @{Ref = "D17"; Val = "2.978.83"},
    @{Ref = "E17"; Val = "  -2.63%  "},
    @{Ref = "D18"; Val = "5.96"},
    @{Ref = "E18"; Val = "  -4.32%  "},
    @{Ref = "D19"; Val = "12.39"},
    @{Ref = "E19"; Val = "  -6.46%  "},
    @{Ref = "D20"; Val = "7.69"},
    @{Ref = "E20"; Val = "  -6.22%  "},
    @{Ref = "D21"; Val = "315.01"},
    @{Ref = "E21"; Val = "  -7.79%  "},
    @{Ref = "D22"; Val = "0.999"},
    @{Ref = "E22"; Val = "  -0.09%  "},
    @{Ref = "D23"; Val = "5.69"},
    @{Ref = "E23"; Val = "  +0.01%  "},
    @{Ref = "D24"; Val = "0.482"},
    @{Ref = "E24"; Val = "  -4.21%  "},
    @{Ref = "D25"; Val = "62.58"},
    @{Ref = "E25"; Val = "  -3.63%  "},
    @{Ref = "D26"; Val = "1.01"},
    @{Ref = "E26"; Val = "  +0.72%  "},
    @{Ref = "D27"; Val = "0.160"},
    @{Ref = "E27"; Val = "  -6.63%  "},
    @{Ref = "D28"; Val = "0.0₃0860"},
    @{Ref = "E28"; Val = "  -12.35%  "},
    @{Ref = "D29"; Val = "6.43"},
    @{Ref = "E29"; Val = "  -7.76%  "},
    @{Ref = "D30"; Val = "6.91"},
    @{Ref = "E30"; Val = "  -7.09%  "},
    @{Ref = "D31"; Val = "1.74"},
    @{Ref = "E31"; Val = "  -6.68%  "},
    @{Ref = "B32"; Val = "EthereumClassic"},
    @{Ref = "C32"; Val = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"},
    @{Ref = "D32"; Val = "19.79"},
    @{Ref = "E32"; Val = "  -6.27%  "},
    @{Ref = "B33"; Val = "Fetch.AI"},
    @{Ref = "C33"; Val = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"},
    @{Ref = "D33"; Val = "1.13"},
    @{Ref = "E33"; Val = "  -9.46%  "},
    @{Ref = "D34"; Val = "153.16"},
    @{Ref = "E34"; Val = "  -2.58%  "},
    @{Ref = "D35"; Val = "4.45"},
    @{Ref = "E35"; Val = "  -7.22%  "},
    @{Ref = "D36"; Val = "5.63"},
    @{Ref = "E36"; Val = "  -6.52%  "},
    @{Ref = "D37"; Val = "1.20"},
    @{Ref = "E37"; Val = "  -9.63%  "},
    @{Ref = "D38"; Val = "23.67"},
    @{Ref = "D39"; Val = "0.0649"},
    @{Ref = "E39"; Val = "  -8.15%  "},
    @{Ref = "D40"; Val = "37.52"},
    @{Ref = "E40"; Val = "  -0.92%  "},
    @{Ref = "D41"; Val = "2.980.80"},
    @{Ref = "E41"; Val = "  -3.92%  "},
    @{Ref = "E42"; Val = "  +0.27%  "},
    @{Ref = "D43"; Val = "0.637"},
    @{Ref = "E43"; Val = "  -4.45%  "},
    @{Ref = "D44"; Val = "3.63"},
    @{Ref = "E44"; Val = "  -7.19%  "},
    @{Ref = "D45"; Val = "2.144.71"},
    @{Ref = "E45"; Val = "  -8.34%  "},
    @{Ref = "D46"; Val = "1.34"},
    @{Ref = "E46"; Val = "  -9.48%  "},
    @{Ref = "D47"; Val = "5.82"},
    @{Ref = "E47"; Val = "  -3.74%  "},
    @{Ref = "D48"; Val = "0.913"},
    @{Ref = "E48"; Val = "  -11.54%  "},
    @{Ref = "D49"; Val = "0.0229"},
    @{Ref = "E49"; Val = "  -6.36%  "},
    @{Ref = "D50"; Val = "18.77"},
    @{Ref = "E50"; Val = "  -7.14%  "},
    @{Ref = "D51"; Val = "0.0853"},
    @{Ref = "E51"; Val = "  -5.27%  "}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    # Force text storage so numeric-looking strings (e.g. "0.160", "23.67")
    # keep their original formatting instead of being parsed as numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Val
    $cell.ClearFormats()
}
